$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Rows(5).Insert()
$ws.Range("A4:F4").Copy()
$ws.Range("A5:F5").PasteSpecial(-4122)
$ws.Range("A5").Value2 = "Câble micro-USB DATA"
$ws.Range("B5").Value2 = "connectique"
$ws.Range("C5").Value2 = "https://www.gotronic.fr/art-cordon-50-cm-rs105-33657.htm"
$ws.Range("D5").Value2 = 1.95
$ws.Range("E5").Value2 = 1
$ws.Range("F5").Formula = "=D5*E5"
$ws.Range("B11").Select()
